# Rewrites the "Exoplanets" essay into "The Everlasting Influence of History" essay,
# matching the target diff (text swaps, a couple of sentence removals, author/email
# swap, and a trailing blank paragraph).
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Finds $phrase (a whole-word Find match) and deletes it together with
# $before characters immediately preceding it and $after characters
# immediately following it (used to also sweep up adjoining punctuation
# that Find can't match on its own, e.g. a lone ".").
function Delete-Phrase($phrase, $before, $after) {
    $rng = $d.Content
    $rng.Find.Execute($phrase) | Out-Null
    $delRange = $d.Range($rng.Start - $before, $rng.End + $after)
    $delRange.Delete() | Out-Null
}

# Finds $phrase and deletes the $count characters immediately following it,
# leaving $phrase itself intact.
function Delete-After($phrase, $count) {
    $rng = $d.Content
    $rng.Find.Execute($phrase) | Out-Null
    $delRange = $d.Range($rng.End, $rng.End + $count)
    $delRange.Delete() | Out-Null
}

# --- Title ---
Replace-Text "Enigmatic Allure of Exoplanets" "The Everlasting Influence of History: A Journey Through Time"

# --- Author name ---
Replace-Text "Dr" "Martin F"
Replace-Text " Amelia Novak" " Davis"

# --- Author email ---
Replace-Text "amelia" "davismf@edu"
Replace-Text "novak@spatiumobservatory" "org"
Delete-After "org" 4   # removes the trailing ".com" (4 chars) after "org"

# --- Body paragraph 1 ---
Replace-Text "In the vast expanse of the cosmos, a realm of wonder and mystery awaits exploration: the world of exoplanets" "History, a vast and intricate tapestry woven from countless moments, serves as a mirror reflecting the evolution of civilizations, the rise and fall of empires, and the indomitable spirit of humanity"
Replace-Text "These celestial bodies, residing beyond our solar system, hold the potential to unveil secrets that could profoundly reshape our understanding of the universe" "As we embark on this journey through time, let us unlock the secrets held within the annals of history, unraveling the mysteries that connect past, present, and future"

# Removes ". With each new discovery, ... scientific knowledge." in full
# (leading "." and trailing "." are single-character runs Find can't match
# alone, so sweep them up via the before/after padding).
Delete-Phrase "With each new discovery, astronomers embark on an exhilarating journey to unravel the complexities of exoplanetary systems, igniting our imaginations and propelling us towards the frontiers of scientific knowledge" 2 1

Replace-Text "Exoplanets present a captivating spectrum of diversity, from gas giants larger than Jupiter to rocky worlds smaller than Earth" "The study of history is a profound undertaking, encompassing diverse civilizations spread across continents and spanning millennia"
Replace-Text "Their existence challenges our notions of planetary formation and habitability, prompting profound questions about the prevalence of life beyond our pale blue dot" "Each civilization, with its unique customs, beliefs, and achievements, adds a vibrant thread to the grand narrative of humanity"
Replace-Text "As we delve deeper into their enigmatic characteristics, we encounter extreme environments that test the limits of our comprehension, revealing phenomena such as scorching hot Jupiters and tidally locked worlds locked in eternal day or night" "From the ancient Egyptians constructing magnificent pyramids to the Renaissance sparking a rebirth of art and science, history showcases the boundless creativity and resilience of humankind"

Replace-Text "The exploration of exoplanets has yielded awe-inspiring insights into the cosmos" "Moreover, history provides invaluable lessons for navigating the complexities of the present and envisioning a more equitable future"
Replace-Text "The discovery of Earth-like exoplanets, potentially capable of supporting liquid water on their surfaces, has fueled our aspirations for finding life beyond Earth" "By examining triumphs and missteps of the past, we glean wisdom that informs our decision-making, fostering empathy and understanding among peoples"
Replace-Text "By scrutinizing the atmospheres of these distant worlds, scientists are unraveling the chemical composition of their atmospheres, searching for telltale signs of biological activity" "Furthermore, history cultivates a sense of continuity, connecting us to those who came before us and reminding us that we are part of an ongoing story"

# Removes ". The analysis of exoplanet transits ... masses and densities." in
# full, but keeps the paragraph's final "." intact.
Delete-Phrase "The analysis of exoplanet transits has provided valuable information about their orbital parameters, allowing us to infer their masses and densities" 2 0

# --- Summary paragraph ---
Replace-Text "The study of exoplanets continues to captivate the imaginations of scientists and enthusiasts alike" "This essay explores the enduring influence of history, highlighting its multifaceted value in shaping our understanding of the world and ourselves"
Replace-Text "As we venture further into the celestial tapestry, unveiling the secrets of these enigmatic worlds, we gain invaluable insights into the fundamental principles governing the formation and evolution of planetary systems, shedding light on our own place in the cosmos" "Through an examination of diverse civilizations and eras, we gain insights into the complexities of human interactions, the evolution of ideas, and the enduring impact of past events on our present circumstances"
Replace-Text "With each new discovery, we edge closer to answering profound questions about the prevalence of life beyond Earth, expanding our understanding of the universe and our position within it" "Ultimately, history serves as an invaluable tool for learning, empathy, and cultivating a sense of interconnectedness, shaping our perspectives and inspiring us to build a better future"

# --- Trailing blank paragraph added after the summary paragraph ---
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter() | Out-Null
